# Apply the changes described by the diff:
#  - workbook.xml: active tab moves from MAY (index 1) to ARDUINO (index 2)
#  - sheet2 (MAY): selection moves from D11 to D14 (and is no longer the active tab)
#  - sheet3 (ARDUINO): new rows 20-21 added (Socket programming unit),
#    dimension grows to A1:H21, view scrolls/selects G18, tab becomes active
#  - sharedStrings.xml: 5 new strings appended for the new content

$wb = $excel.ActiveWorkbook

$wsMay = $wb.Worksheets.Item("MAY")
$wsArduino = $wb.Worksheets.Item("ARDUINO")

# MAY sheet keeps its own remembered selection (D14) but stops being the
# tab that is active/selected once we switch away from it below.
$wsMay.Range("D14").Select() | Out-Null

# Make ARDUINO the active sheet (drives workbook.xml activeTab + this
# sheet's tabSelected attribute).
$wsArduino.Select() | Out-Null

# Add the new "[Unit 9] Socket programming" entries in rows 20-21.
# Values/hyperlinks are entered in an order that matches the order new
# shared strings were appended in the source workbook.
$wsArduino.Range("H20").Value = "[Unit 9] Socket programming"
$wsArduino.Hyperlinks.Add($wsArduino.Range("G20"), "https://www.youtube.com/watch?v=-Fs6wAV7tEw") | Out-Null

$wsArduino.Range("F20").Value = "Socket"

$wsArduino.Hyperlinks.Add($wsArduino.Range("G21"), "https://www.youtube.com/watch?v=LWdynDo5jqo") | Out-Null
$wsArduino.Range("H21").Value = "[Đồ án mạng máy tính] Share 1: Làm quen với Python socket"

# Match the existing hyperlink cell formatting used elsewhere in column G.
$wsArduino.Range("G20").Style = "Hyperlink"
$wsArduino.Range("G21").Style = "Hyperlink"

# Final selection on the ARDUINO sheet.
$wsArduino.Range("G18").Select() | Out-Null
